$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-15 19:19:05"
$ws.Range("E3").Value = "2026-02-15 19:19:08"
$ws.Range("I3").Value = "1.0 mm"
$ws.Range("O3").Value = "-5.3 °C"
$ws.Range("E4").Value = "2026-02-15 19:19:11"
$ws.Range("H4").Value = "'70%"
$ws.Range("E5").Value = "2026-02-15 19:19:14"
$ws.Range("H5").Value = "'93%"
$ws.Range("I5").Value = "4.0 mm"
$ws.Range("L5").Value = "34.9 km/h - 322º 18:38 TU"
$ws.Range("O5").Value = "-4.7 °C"
$ws.Range("E6").Value = "2026-02-15 19:19:17"
$ws.Range("O6").Value = "8.6 °C"
$ws.Range("E7").Value = "2026-02-15 19:19:20"
$ws.Range("O7").Value = "11.6 °C"
$ws.Range("E8").Value = "2026-02-15 19:19:23"
$ws.Range("O8").Value = "8.1 °C"
$ws.Range("E9").Value = "2026-02-15 19:19:25"
$ws.Range("H9").Value = "'49%"
$ws.Range("E10").Value = "2026-02-15 19:19:28"
$ws.Range("E11").Value = "2026-02-15 19:19:31"
$ws.Range("H11").Value = "'41%"
$ws.Range("E12").Value = "2026-02-15 19:19:33"
$ws.Range("H12").Value = "'54%"
$ws.Range("E13").Value = "2026-02-15 19:19:36"
$ws.Range("E14").Value = "2026-02-15 19:19:39"
$ws.Range("H14").Value = "'58%"
$ws.Range("O14").Value = "10.6 °C"
$ws.Range("E15").Value = "2026-02-15 19:19:41"
$ws.Range("H15").Value = "'49%"
$ws.Range("E16").Value = "2026-02-15 19:19:44"
$ws.Range("E17").Value = "2026-02-15 19:19:47"
$ws.Range("E18").Value = "2026-02-15 19:19:49"
$ws.Range("H18").Value = "'71%"
$ws.Range("O18").Value = "7.4 °C"
$ws.Range("E19").Value = "2026-02-15 19:19:52"
$ws.Range("O19").Value = "3.3 °C"
$ws.Range("E20").Value = "2026-02-15 19:19:55"
$ws.Range("O20").Value = "-2.8 °C"
$ws.Range("E21").Value = "2026-02-15 19:19:58"
$ws.Range("E22").Value = "2026-02-15 19:20:01"
$ws.Range("E23").Value = "2026-02-15 19:20:04"
$ws.Range("I23").Value = "1.9 mm"
$ws.Range("O23").Value = "-3.7 °C"
$ws.Range("E24").Value = "2026-02-15 19:20:07"
$ws.Range("O24").Value = "8.8 °C"
$ws.Range("E25").Value = "2026-02-15 19:20:09"
$ws.Range("O25").Value = "-1.7 °C"
$ws.Range("E26").Value = "2026-02-15 19:20:12"
$ws.Range("E27").Value = "2026-02-15 19:20:14"
$ws.Range("H27").Value = "'48%"
$ws.Range("O27").Value = "-0.1 °C"
$ws.Range("E28").Value = "2026-02-15 19:20:17"
$ws.Range("H28").Value = "'57%"
$ws.Range("J28").Value = "1015.6 hPa"
$ws.Range("E29").Value = "2026-02-15 19:20:20"
$ws.Range("H29").Value = "'56%"
$ws.Range("O29").Value = "10.1 °C"
$ws.Range("E30").Value = "2026-02-15 19:20:22"
$ws.Range("H30").Value = "'53%"
$ws.Range("E31").Value = "2026-02-15 19:20:25"
$ws.Range("O31").Value = "9.9 °C"
$ws.Range("E32").Value = "2026-02-15 19:20:28"
$ws.Range("H32").Value = "'83%"
$ws.Range("O32").Value = "3.6 °C"
$ws.Range("E33").Value = "2026-02-15 19:20:32"
$ws.Range("H33").Value = "'41%"
$ws.Range("J33").Value = "1015.2 hPa"
$ws.Range("O33").Value = "5.7 °C"
$ws.Range("E34").Value = "2026-02-15 19:20:34"
$ws.Range("M34").Value = "4.3 °C 18:59 TU"
$ws.Range("E35").Value = "2026-02-15 19:20:37"
$ws.Range("J35").Value = "1019.5 hPa"
$ws.Range("E36").Value = "2026-02-15 19:20:40"
$ws.Range("H36").Value = "'47%"
$ws.Range("E37").Value = "2026-02-15 19:20:43"
$ws.Range("E38").Value = "2026-02-15 19:20:46"
$ws.Range("E39").Value = "2026-02-15 19:20:48"
$ws.Range("H39").Value = "'58%"
$ws.Range("M39").Value = "1.3 °C 18:51 TU"
$ws.Range("O39").Value = "-3.0 °C"
$ws.Range("E40").Value = "2026-02-15 19:20:51"
$ws.Range("H40").Value = "'35%"
$ws.Range("J40").Value = "1016.0 hPa"
$ws.Range("O40").Value = "9.0 °C"
$ws.Range("E41").Value = "2026-02-15 19:20:54"
$ws.Range("O41").Value = "12.4 °C"
$ws.Range("E42").Value = "2026-02-15 19:20:57"
$ws.Range("E43").Value = "2026-02-15 19:21:00"
$ws.Range("O43").Value = "6.2 °C"
$ws.Range("E44").Value = "2026-02-15 19:21:02"
$ws.Range("I44").Value = "1.9 mm"
$ws.Range("E45").Value = "2026-02-15 19:21:05"
$ws.Range("I45").Value = "0.7 mm"
$ws.Range("O45").Value = "0.9 °C"
$ws.Range("E46").Value = "2026-02-15 19:21:08"
$ws.Range("H46").Value = "'52%"
$ws.Range("O46").Value = "11.7 °C"
